$wb = $excel.ActiveWorkbook

# ---------- Sheet "2006" (sheet1) ----------
$ws1 = $wb.Worksheets.Item(1)

# Insert new column D (time) before the existing data columns.
$ws1.Columns.Item(4).Insert()
$ws1.Columns.Item(4).ColumnWidth = 5.140625

$ws1.Range("D1").Value = "time"
$ws1.Range("D2:D25").Value = 0

# Insert 24 new rows after row 25 for the 2006 (12-month) follow-up data.
$ws1.Rows.Item(26).Resize(24).Insert()

$ws1.Range("A26").Value = "Information Processing Speed"
$ws1.Range("B26").Value = "Stroop color word test mean card 1+2 b"
$ws1.Range("D26").Value = 12
$ws1.Range("E26").Value = "49.0 (8.9)"
$ws1.Range("F26").Value = "49.7 (7.7)"
$ws1.Range("G26").Value = "47.8 (9.1)"
$ws1.Range("H26").Value = "44.8 (6.6)"

$ws1.Range("A27").Value = "Information Processing Speed"
$ws1.Range("B27").Value = "WAIS Digit Symbol c"
$ws1.Range("D27").Value = 12
$ws1.Range("E27").Value = "59.0 (11.9)"
$ws1.Range("F27").Value = "57.2 (8.2)"
$ws1.Range("G27").Value = "58.7 (12.3)"
$ws1.Range("H27").Value = "62.6 (10.3)"

$ws1.Range("A28").Value = "Attention"
$ws1.Range("B28").Value = "Trailmaking A b"
$ws1.Range("D28").Value = 12
$ws1.Range("E28").Value = "28.0 (8.7)"
$ws1.Range("F28").Value = "28.2 (8.5)"
$ws1.Range("G28").Value = "30.8 (9.3)"
$ws1.Range("H28").Value = "28.6 (10.5)"

$ws1.Range("A29").Value = "Attention"
$ws1.Range("B29").Value = "Eriksen congruent b"
$ws1.Range("C29").Value = "This is not a standard clinical NP test - we may want to take out"
$ws1.Range("D29").Value = 12
$ws1.Range("E29").Value = "559.1 (51.9)"
$ws1.Range("F29").Value = "571.9 (57.5)"
$ws1.Range("G29").Value = "570.1 (60.3)"
$ws1.Range("H29").Value = "551,0 (59.2)"

$ws1.Range("A30").Value = "Attention"
$ws1.Range("B30").Value = "Eriksen perceptual b"
$ws1.Range("D30").Value = 12
$ws1.Range("E30").Value = "563.5 (50.9)"
$ws1.Range("F30").Value = "581.0 (51.9)"
$ws1.Range("G30").Value = "575.7 (58.9)"
$ws1.Range("H30").Value = "558.9 (59.0)"

$ws1.Range("A31").Value = "Attention"
$ws1.Range("B31").Value = "Eriksen response conflict b"
$ws1.Range("D31").Value = 12
$ws1.Range("E31").Value = "594.0 (50.5)"
$ws1.Range("F31").Value = "609.5 (48.4)"
$ws1.Range("G31").Value = "601.9 (56.3)"
$ws1.Range("H31").Value = "580.4 (56.1)"

$ws1.Range("A32").Value = "Executive Functions"
$ws1.Range("B32").Value = "Stroop color word test card 4 b"
$ws1.Range("D32").Value = 12
$ws1.Range("E32").Value = "104.6 (27.0)"
$ws1.Range("F32").Value = "100.0 (19.8)"
$ws1.Range("G32").Value = "107.4 (26.4)"
$ws1.Range("H32").Value = "91.2 (15.6)"

$ws1.Range("A33").Value = "Executive Functions"
$ws1.Range("B33").Value = "Stroop color word test interference b"
$ws1.Range("D33").Value = 12
$ws1.Range("E33").Value = "1.6 (0.2)"
$ws1.Range("F33").Value = "1.6 (0.2)"
$ws1.Range("G33").Value = "1.6 (0.2)"
$ws1.Range("H33").Value = "1.6 (0.2)"

$ws1.Range("A34").Value = "Executive Functions"
$ws1.Range("B34").Value = "Trailmaking B b"
$ws1.Range("D34").Value = 12
$ws1.Range("E34").Value = "66.8 (23.1)"
$ws1.Range("F34").Value = "60.1 (18.1)"
$ws1.Range("G34").Value = "73.7 (28.1)"
$ws1.Range("H34").Value = "62.1 (22.3)"

$ws1.Range("A35").Value = "??"
$ws1.Range("B35").Value = "AFM task stimulus identification b"
$ws1.Range("C35").Value = "This task seems to be an experimental measure - not sure if it should be included"
$ws1.Range("D35").Value = 12
$ws1.Range("E35").Value = "113.6 (29.8)"
$ws1.Range("F35").Value = "113.6 (29.1)"
$ws1.Range("G35").Value = "106.1 (39.8)"
$ws1.Range("H35").Value = "112.2 (30.1)"

$ws1.Range("A36").Value = "??"
$ws1.Range("B36").Value = "AFM task central response decision b"
$ws1.Range("D36").Value = 12
$ws1.Range("E36").Value = "58.4 (35.6)"
$ws1.Range("F36").Value = "55.6 (36.7)"
$ws1.Range("G36").Value = "62.5 (49.2)"
$ws1.Range("H36").Value = "61.1 (46.3)"

$ws1.Range("A37").Value = "??"
$ws1.Range("B37").Value = "AFM task time response preparation b"
$ws1.Range("D37").Value = 12
$ws1.Range("E37").Value = "11.2 (28.1)"
$ws1.Range("F37").Value = "5.1 (32.8)"
$ws1.Range("G37").Value = "14.1 (35.6)"
$ws1.Range("H37").Value = "12.8 (33.4)"

$ws1.Range("A38").Value = "??"
$ws1.Range("B38").Value = "AFM reaction time corrected for errors b"
$ws1.Range("D38").Value = 12
$ws1.Range("E38").Value = "616.4 (80.7)"
$ws1.Range("F38").Value = "650.7 (96.8)"
$ws1.Range("G38").Value = "651.4 (107.9)"
$ws1.Range("H38").Value = "599.9 (90.4)"

$ws1.Range("A39").Value = "Verbal Memory"
$ws1.Range("B39").Value = "CVLT recall c"
$ws1.Range("D39").Value = 12
$ws1.Range("E39").Value = "59.2 (9.3)"
$ws1.Range("F39").Value = "57.2 (9.7)"
$ws1.Range("G39").Value = "57.8 (8.7)"
$ws1.Range("H39").Value = "59.9 (8.4)"

$ws1.Range("A40").Value = "Verbal Memory"
$ws1.Range("B40").Value = "CVLT delayed recall c"
$ws1.Range("D40").Value = 12
$ws1.Range("E40").Value = "13.8 ( 2.4)"
$ws1.Range("F40").Value = "13.6 ( 2.5)"
$ws1.Range("G40").Value = "13.2 (2.7)"
$ws1.Range("H40").Value = "14.1 (2.2)"

$ws1.Range("A41").Value = "Verbal Memory"
$ws1.Range("B41").Value = "CVLT recognition c"
$ws1.Range("D41").Value = 12
$ws1.Range("E41").Value = "15.6 (0.7)"
$ws1.Range("F41").Value = "15.5 ( 0.7)"
$ws1.Range("G41").Value = "15.2 (1.1)"
$ws1.Range("H41").Value = "15.5 (0.9)"

$ws1.Range("A42").Value = "Visual Memory"
$ws1.Range("B42").Value = "WMS immediate recall c "
$ws1.Range("C42").Value = "According to de Ruiter this is Visual Reproduction I"
$ws1.Range("D42").Value = 12
$ws1.Range("E42").Value = "34.6 (4.3) "
$ws1.Range("F42").Value = "33.8 (3.5) "
$ws1.Range("G42").Value = " 33.2 (4.7)"
$ws1.Range("H42").Value = "35.3 (3.6)"

$ws1.Range("A43").Value = "Visual Memory"
$ws1.Range("B43").Value = "WMS delayed recall c"
$ws1.Range("C43").Value = "According to de Ruiter this is Visual Reproduction II"
$ws1.Range("D43").Value = 12
$ws1.Range("E43").Value = "31.7 (6.2)"
$ws1.Range("F43").Value = "32.1 (5.6)"
$ws1.Range("G43").Value = "31.2 (5.7)"
$ws1.Range("H43").Value = "33.6 (4.6)"

$ws1.Range("A44").Value = "??"
$ws1.Range("B44").Value = "Memory update 6 c"
$ws1.Range("C44").Value = "I don't know what this test is"
$ws1.Range("D44").Value = 12
$ws1.Range("E44").Value = "86.8 (11.5)"
$ws1.Range("F44").Value = "88.2 (9.6)"
$ws1.Range("G44").Value = "84.6 (12.3)"
$ws1.Range("H44").Value = "89.2 (11.1)"

$ws1.Range("A45").Value = "??"
$ws1.Range("B45").Value = "Memory update 9 c"
$ws1.Range("D45").Value = 12
$ws1.Range("E45").Value = "76.6 (11.8)"
$ws1.Range("F45").Value = "78.6 (10.5)"
$ws1.Range("G45").Value = "76.0 (13.5)"
$ws1.Range("H45").Value = "80.8 (11.5)"

$ws1.Range("A46").Value = "??"
$ws1.Range("B46").Value = "Memory update 12 c"
$ws1.Range("D46").Value = 12
$ws1.Range("E46").Value = "73.4 (12.1)"
$ws1.Range("F46").Value = "74.5 (13.2)"
$ws1.Range("G46").Value = "68.7 (13.8)"
$ws1.Range("H46").Value = "74.9 (10.3)"

$ws1.Range("A47").Value = "Language"
$ws1.Range("B47").Value = "Word fluency c"
$ws1.Range("D47").Value = 12
$ws1.Range("E47").Value = "23.0 (5.6)"
$ws1.Range("F47").Value = "21.9 ( 4.6)"
$ws1.Range("G47").Value = "22.3 (4.8)"
$ws1.Range("H47").Value = "23.0 (5.0)"

$ws1.Range("A48").Value = "Motor Speed"
$ws1.Range("B48").Value = "Fepsy Finger Tapping (dominant) c"
$ws1.Range("D48").Value = 12
$ws1.Range("E48").Value = "57.8 (7.5)"
$ws1.Range("F48").Value = "56.9 (7.7)"
$ws1.Range("G48").Value = "58.6 (7.8)"
$ws1.Range("H48").Value = "56.2 (6.7)"

$ws1.Range("A49").Value = "Motor Speed"
$ws1.Range("B49").Value = "Fepsy Finger Tapping (non dominant) c"
$ws1.Range("D49").Value = 12
$ws1.Range("E49").Value = "51.0 (6.4)"
$ws1.Range("F49").Value = "51.2 (9.5)"
$ws1.Range("G49").Value = "52.0 (8.3)"
$ws1.Range("H49").Value = "51.4 (7.3)"

$ws1.Range("A26:A49").HorizontalAlignment = -4108
$ws1.Range("C3").Select()

# ---------- Sheet "2015" (sheet2) ----------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Columns.Item(4).Insert()
$ws2.Range("D1").Value = "time"
$ws2.Range("D2").Formula = "=11.5*12"
$ws2.Range("D3:D16").Formula = "=11.5*12"
